$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("2:2").Insert()
$ws.Rows("2:2").RowHeight = 22.05
